$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Harpreet Brar"

# Insert a new first column (shifts teamName..result from A:L to B:M)
$ws.Columns.Item(1).Insert()

# Insert a new row above the existing data row so the existing record
# (currently row 2, now holding B2:M2 = Punjab Kings/Harpreet Brar/.../Mumbai won...)
# becomes row 3, and we can place the new "37th" record into row 2.
$ws.Rows.Item(2).Insert()

# All rows of data, in final column order:
# matchNo, teamName, batterName, states, runs, balls, fours, sixes, sr, opponentTeamName, venue, date, result
$allRows = @(
    @("matchNo", "teamName", "batterName", "states", "runs", "balls", "fours", "sixes", "sr", "opponentTeamName", "venue", "date", "result"),
    @("37th", "Punjab Kings", "Harpreet Brar", "", "18", "18", "1", "0", "100.00", "Sunrisers Hyderabad", "Sharjah", "September 25", "Punjab Kings won by 5 runs"),
    @("42nd", "Punjab Kings", "Harpreet Brar", "", "14", "19", "0", "0", "73.68", "Mumbai Indians", "Abu Dhabi", "September 28", "Mumbai won by 6 wickets (with 6 balls remaining)"),
    @("48th", "Punjab Kings", "Harpreet Brar", "", "3", "2", "0", "0", "150.00", "Royal Challengers Bangalore", "Sharjah", "October 03", "RCB won by 6 runs"),
    @("26th", "Punjab Kings", "Harpreet Brar", "", "25", "17", "1", "2", "147.05", "Royal Challengers Bangalore", "Ahmedabad", "April 30", "Punjab Kings won by 34 runs"),
    @("29th", "Punjab Kings", "Harpreet Brar", "", "4", "2", "1", "0", "200.00", "Delhi Capitals", "Ahmedabad", "May 02", "Capitals won by 7 wickets (with 14 balls remaining)")
)

# Columns E:I (runs, balls, fours, sixes, sr) hold numeric-looking text
# ("18", "100.00", ...). Format those data cells (rows 2-6) as Text BEFORE
# writing so Excel keeps them as strings instead of auto-converting to
# numbers - matches the source data's numberStoredAsText convention.
$ws.Range("E2:I6").NumberFormat = "@"

for ($i = 0; $i -lt $allRows.Count; $i++) {
    $r = $i + 1
    $rowValues = $allRows[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}
